$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.915.47"
$ws.Range("E2").Value = "'  -0.10%  "
$ws.Range("D3").Value = "'1.813.31"
$ws.Range("E3").Value = "'  +0.27%  "
$ws.Range("E4").Value = "'  +0.16%  "
$ws.Range("D5").Value = "'309.14"
$ws.Range("E5").Value = "'  -0.47%  "
$ws.Range("E6").Value = "'  +0.06%  "
$ws.Range("E7").Value = "'  +0.09%  "
$ws.Range("D8").Value = "'0.3657"
$ws.Range("E8").Value = "'  -1.65%  "
$ws.Range("D9").Value = "'0.07361"
$ws.Range("E9").Value = "'  -0.15%  "
$ws.Range("D10").Value = "'0.8679"
$ws.Range("E10").Value = "'  -0.86%  "
$ws.Range("D11").Value = "'20.21"
$ws.Range("E11").Value = "'  -1.25%  "
$ws.Range("D12").Value = "'1.835.74"
$ws.Range("E12").Value = "'  +1.73%  "
$ws.Range("E13").Value = "'  -0.20%  "
$ws.Range("D14").Value = "'0.07099"
$ws.Range("E14").Value = "'  +0.56%  "
$ws.Range("D15").Value = "'6.497"
$ws.Range("E15").Value = "'  -0.48%  "
$ws.Range("D16").Value = "'90.94"
$ws.Range("E16").Value = "'  -1.69%  "
$ws.Range("D17").Value = "'1.004"
$ws.Range("D18").Value = "'0.000008693"
$ws.Range("E18").Value = "'  -0.25%  "
$ws.Range("E19").Value = "'  +0.10%  "
$ws.Range("D20").Value = "'14.62"
$ws.Range("E20").Value = "'  -0.67%  "
$ws.Range("D21").Value = "'26.931.49"
$ws.Range("E21").Value = "'  -0.06%  "
$ws.Range("D22").Value = "'5.282"
$ws.Range("E22").Value = "'  -0.49%  "
$ws.Range("E23").Value = "'  -0.63%  "
$ws.Range("D24").Value = "'2.047.07"
$ws.Range("E24").Value = "'  +0.41%  "
$ws.Range("D25").Value = "'1.892"
$ws.Range("E25").Value = "'  -0.85%  "
$ws.Range("D26").Value = "'150.78"
$ws.Range("E26").Value = "'  -0.60%  "
$ws.Range("E27").Value = "'  -0.59%  "
$ws.Range("D28").Value = "'2.120"
$ws.Range("E28").Value = "'  -1.41%  "
$ws.Range("D29").Value = "'5.243"
$ws.Range("E29").Value = "'  -0.94%  "
$ws.Range("D30").Value = "'115.29"
$ws.Range("E30").Value = "'  -0.57%  "
$ws.Range("E31").Value = "'  -0.26%  "
$ws.Range("D32").Value = "'0.7526"
$ws.Range("E32").Value = "'  -0.58%  "
$ws.Range("D33").Value = "'1.160"
$ws.Range("E33").Value = "'  +0.29%  "
$ws.Range("D34").Value = "'4.476"
$ws.Range("E34").Value = "'  +0.34%  "
$ws.Range("E35").Value = "'  -1.05%  "
$ws.Range("E36").Value = "'  +0.12%  "
$ws.Range("D37").Value = "'1.094"
$ws.Range("E37").Value = "'  -1.02%  "
$ws.Range("D38").Value = "'0.05284"
$ws.Range("E38").Value = "'  +0.72%  "
$ws.Range("D39").Value = "'0.01943"
$ws.Range("E39").Value = "'  -1.63%  "
$ws.Range("D40").Value = "'2.979"
$ws.Range("E40").Value = "'  +1.96%  "
$ws.Range("D41").Value = "'7.231"
$ws.Range("E41").Value = "'  +0.20%  "
$ws.Range("D42").Value = "'0.5285"
$ws.Range("E42").Value = "'  -0.70%  "
$ws.Range("D43").Value = "'2.287"
$ws.Range("E43").Value = "'  -5.81%  "
$ws.Range("D44").Value = "'0.1651"
$ws.Range("E44").Value = "'  -0.93%  "
$ws.Range("D45").Value = "'8.396"
$ws.Range("E45").Value = "'  -1.41%  "
$ws.Range("D46").Value = "'0.4848"
$ws.Range("E46").Value = "'  -3.11%  "
$ws.Range("D47").Value = "'10.46"
$ws.Range("E47").Value = "'  +1.17%  "
$ws.Range("E48").Value = "'  +0.07%  "
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").Value = "'103.16"
$ws.Range("E49").Value = "'  -0.81%  "
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").Value = "'1.658"
$ws.Range("E50").Value = "'  -1.17%  "
$ws.Range("E51").Value = "'  -0.09%  "
